$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.375.43"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").Value = "1.653.22"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.51%  "

$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").Value = "1.887.55"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").Value = "1.652.11"
$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.570"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "27.386.50"
$ws.Range("E17").Value = "  -1.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.38%  "

$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("E23").Value = "  +1.49%  "

$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.85%  "

$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("D33").Value = "1.459.60"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.908"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "

$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.31%  "

$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "1.795.54"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "
